$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking price strings (e.g. "1.00", "2.80")
# are stored verbatim as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.889.49"
$ws.Range("E2").Value = "  -1.43%  "
$ws.Range("D3").Value = "2.223.66"
$ws.Range("E3").Value = "  -0.72%  "
$ws.Range("E4").Value = "  -1.90%  "
$ws.Range("D5").Value = "298.89"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").Value = "90.41"
$ws.Range("E6").Value = "  -3.90%  "
$ws.Range("D7").Value = "0.553"
$ws.Range("E7").Value = "  -3.18%  "
$ws.Range("E8").Value = "  -0.65%  "
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  -6.06%  "
$ws.Range("D10").Value = "32.94"
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("E11").Value = "  -3.52%  "
$ws.Range("D12").Value = "6.94"
$ws.Range("E12").Value = "  -3.75%  "
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").Value = "2.563.05"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "2.225.85"
$ws.Range("E15").Value = "  -3.41%  "
$ws.Range("D16").Value = "13.43"
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "0.775"
$ws.Range("E17").Value = "  -7.02%  "
$ws.Range("D18").Value = "43.768.66"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("D19").Value = "0.0₃0900"
$ws.Range("E19").Value = "  -5.73%  "
$ws.Range("D20").Value = "11.50"
$ws.Range("E20").Value = "  -3.89%  "
$ws.Range("D21").Value = "5.93"
$ws.Range("E21").Value = "  -6.47%  "
$ws.Range("D22").Value = "64.43"
$ws.Range("E22").Value = "  -2.08%  "
$ws.Range("D23").Value = "235.99"
$ws.Range("E23").Value = "  -0.82%  "
$ws.Range("D24").Value = "2.80"
$ws.Range("E24").Value = "  -5.34%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  -0.12%  "
$ws.Range("D26").Value = "1.87"
$ws.Range("E26").Value = "  -5.42%  "
$ws.Range("D27").Value = "38.16"
$ws.Range("E27").Value = "  +1.63%  "
$ws.Range("D28").Value = "2.21"
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("D29").Value = "9.33"
$ws.Range("E29").Value = "  -4.69%  "
$ws.Range("D30").Value = "152.92"
$ws.Range("E30").Value = "  +0.22%  "
$ws.Range("D31").Value = "19.21"
$ws.Range("E31").Value = "  -3.79%  "
$ws.Range("D32").Value = "5.41"
$ws.Range("E32").Value = "  -8.75%  "
$ws.Range("E33").Value = "  -5.43%  "
$ws.Range("D34").Value = "2.50"
$ws.Range("E34").Value = "  -5.36%  "
$ws.Range("D35").Value = "0.116"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").Value = "2.83"
$ws.Range("E36").Value = "  -8.21%  "
$ws.Range("E37").Value = "  -6.40%  "
$ws.Range("E38").Value = "  -5.43%  "
$ws.Range("D39").Value = "0.0298"
$ws.Range("E39").Value = "  -0.65%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  -6.02%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "3.62"
$ws.Range("E41").Value = "  -4.15%  "
$ws.Range("D42").Value = "13.18"
$ws.Range("E42").Value = "  -8.61%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.93%  "
$ws.Range("D44").Value = "1.831.76"
$ws.Range("E44").Value = "  +3.95%  "
$ws.Range("D45").Value = "1.78"
$ws.Range("E45").Value = "  +12.39%  "
$ws.Range("E46").Value = "  -6.06%  "
$ws.Range("B47").Value = "ordi"
$ws.Range("C47").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D47").Value = "67.16"
$ws.Range("E47").Value = "  -3.51%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "94.32"
$ws.Range("E48").Value = "  -4.53%  "
$ws.Range("D49").Value = "72.87"
$ws.Range("E49").Value = "  -8.14%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "13.80"
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").Value = "7.69"
$ws.Range("E51").Value = "  -4.78%  "
